# Response to Reviewers.docx -- "LP very minor edits"
#
# 1) The last live edit in the letter happened in the opening salutation
#    ("Dear Editors and Reviewers,"); Word always drops its "_GoBack"
#    bookmark (last-edit marker) at the point of the most recent edit,
#    removing it from wherever it previously sat. We reproduce that by
#    relocating the bookmark from its old spot (inside the Reviewer 2
#    paragraph) to just after "Dear Editor" in the salutation.
# 2) Track changes are cleaned up (the lone insertion mark by
#    Vigers, Timothy is accepted).
# 3) "hopefully" -> "now" in the figure-legend-rendering sentence.
# 4) The closing signature line is reordered from
#    "Laura Pyle and Tim Vigers" to "Tim Vigers and Laura Pyle".

$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark to the salutation -------------------
# (Bookmark names are unique; re-adding "_GoBack" at the new location
# automatically removes it from its old spot, exactly like Word's
# last-edit-location tracking.)
$salutation = $d.Content
$salutation.Find.Execute("Dear Editor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$salutation.Collapse(0)
$d.Bookmarks.Add("_GoBack", $salutation)

# --- 2. Accept all tracked changes (removes the stray <w:ins/> marker) -
$d.Revisions.AcceptAll()

# --- 3. "hopefully" -> "now" --------------------------------------------
$d.Content.Find.Execute("Figure legends should hopefully render correctly", $true, $false, $false, $false, $false, $true, 1, $false, "Figure legends should now render correctly", 2)

# --- 4. Reorder the sign-off names --------------------------------------
$d.Content.Find.Execute("Laura Pyle and Tim Vigers", $true, $false, $false, $false, $false, $true, 1, $false, "Tim Vigers and Laura Pyle", 2)
